# RSII - Lesson 04: "Unos 5-tog i 6tog prdavana"
#
# 1) Bump the cached "datetimeFigureOut" field text from 3/29/2019 to
#    3/30/2019 on the slide master, every slide layout and the notes
#    master (13 placeholders total).
# 2) Fix a typo on slide 5 ("Moview" -> "Movie") inside the C# code
#    sample, merging the three runs that made up "$"Moview (" into a
#    single "$"Movie (" run.
# 3) Fix a typo on slide 9 ("objecta" -> "objekta"), merging the three
#    runs that made up "novog objecta sa ovim " into a single
#    "novog objekta sa ovim " run.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes, $newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "3/30/2019"

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes $newDate

# Every slide layout
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes $newDate
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes $newDate

# --- Slide 5: "$"Moview (" -> "$"Movie (" -------------------------------
$slide5 = $p.Slides.Item(5)
$rect4 = $slide5.Shapes.Item(2)
$tr5 = $rect4.TextFrame.TextRange
$full5 = $tr5.Text
$old5 = '$"Moview ('
$new5 = '$"Movie ('
$idx5 = $full5.IndexOf($old5)
if ($idx5 -ge 0) {
    $sub5 = $tr5.Characters($idx5 + 1, $old5.Length)
    $sub5.Text = $new5
}

# --- Slide 9: "novog objecta sa ovim " -> "novog objekta sa ovim " ------
$slide9 = $p.Slides.Item(9)
$contentPh = $slide9.Shapes.Item(2)
$tr9 = $contentPh.TextFrame.TextRange
$full9 = $tr9.Text
$old9 = 'Pri generiranju novog objecta sa ovim '
$new9 = 'Pri generiranju novog objekta sa ovim '
$idx9 = $full9.IndexOf($old9)
if ($idx9 -ge 0) {
    $sub9 = $tr9.Characters($idx9 + 1, $old9.Length)
    $sub9.Text = $new9
}
